$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AgTests (F) and AgPosit (G) values for rows 303-361 per updated source data
$ws.Range("F303").Value = 19194
$ws.Range("G303").Value = 1226
$ws.Range("F304").Value = 12104
$ws.Range("G304").Value = 1042
$ws.Range("F305").Value = 6722
$ws.Range("G305").Value = 514
$ws.Range("F306").Value = 149810
$ws.Range("G306").Value = 15227
$ws.Range("F307").Value = 151488
$ws.Range("G307").Value = 12793
$ws.Range("F308").Value = 30887
$ws.Range("G308").Value = 2094
$ws.Range("F309").Value = 155738
$ws.Range("G309").Value = 11038
$ws.Range("F310").Value = 158451
$ws.Range("G310").Value = 8132
$ws.Range("F311").Value = 122965
$ws.Range("G311").Value = 3842
$ws.Range("F312").Value = 56240
$ws.Range("G312").Value = 1843
$ws.Range("F313").Value = 151115
$ws.Range("G313").Value = 6899
$ws.Range("F314").Value = 128462
$ws.Range("G314").Value = 6288
$ws.Range("F315").Value = 112548
$ws.Range("G315").Value = 5245
$ws.Range("F316").Value = 101457
$ws.Range("G316").Value = 4585
$ws.Range("F317").Value = 127448
$ws.Range("G317").Value = 4341
$ws.Range("F318").Value = 97888
$ws.Range("G318").Value = 2257
$ws.Range("F319").Value = 82644
$ws.Range("G319").Value = 3256
$ws.Range("F320").Value = 143064
$ws.Range("G320").Value = 6573
$ws.Range("F321").Value = 178489
$ws.Range("G321").Value = 5291
$ws.Range("F322").Value = 219170
$ws.Range("G322").Value = 4663
$ws.Range("F323").Value = 433640
$ws.Range("G323").Value = 6216
$ws.Range("F324").Value = 481577
$ws.Range("G324").Value = 5594
$ws.Range("F325").Value = 1531885
$ws.Range("G325").Value = 12917
$ws.Range("F326").Value = 839091
$ws.Range("G326").Value = 7659
$ws.Range("F327").Value = 450134
$ws.Range("G327").Value = 5439
$ws.Range("F328").Value = 361481
$ws.Range("G328").Value = 5331
$ws.Range("F329").Value = 165954
$ws.Range("G329").Value = 3516
$ws.Range("F330").Value = 145016
$ws.Range("G330").Value = 4156
$ws.Range("F331").Value = 309380
$ws.Range("G331").Value = 5398
$ws.Range("F332").Value = 914201
$ws.Range("G332").Value = 9070
$ws.Range("F333").Value = 543422
$ws.Range("G333").Value = 5890
$ws.Range("F334").Value = 196605
$ws.Range("G334").Value = 3495
$ws.Range("F335").Value = 130637
$ws.Range("G335").Value = 2995
$ws.Range("F336").Value = 102457
$ws.Range("G336").Value = 3342
$ws.Range("F337").Value = 103305
$ws.Range("G337").Value = 2879
$ws.Range("F338").Value = 227144
$ws.Range("G338").Value = 3177
$ws.Range("F339").Value = 661348
$ws.Range("G339").Value = 5502
$ws.Range("F340").Value = 385238
$ws.Range("G340").Value = 3309
$ws.Range("F341").Value = 291263
$ws.Range("G341").Value = 3652
$ws.Range("F342").Value = 178143
$ws.Range("G342").Value = 3029
$ws.Range("F343").Value = 132807
$ws.Range("G343").Value = 2964
$ws.Range("F344").Value = 135000
$ws.Range("G344").Value = 2474
$ws.Range("F345").Value = 291585
$ws.Range("G345").Value = 3314
$ws.Range("F346").Value = 675177
$ws.Range("G346").Value = 4832
$ws.Range("F347").Value = 343715
$ws.Range("G347").Value = 2920
$ws.Range("F348").Value = 231830
$ws.Range("G348").Value = 3232
$ws.Range("F349").Value = 158692
$ws.Range("G349").Value = 2747
$ws.Range("F350").Value = 126706
$ws.Range("G350").Value = 2777
$ws.Range("F351").Value = 150659
$ws.Range("G351").Value = 2821
$ws.Range("F352").Value = 307430
$ws.Range("G352").Value = 3541
$ws.Range("F353").Value = 725102
$ws.Range("G353").Value = 5296
$ws.Range("F354").Value = 312740
$ws.Range("G354").Value = 2866
$ws.Range("F355").Value = 221562
$ws.Range("G355").Value = 3439
$ws.Range("F356").Value = 159526
$ws.Range("G356").Value = 2868
$ws.Range("F357").Value = 138102
$ws.Range("G357").Value = 3021
$ws.Range("F358").Value = 157953
$ws.Range("G358").Value = 2600
$ws.Range("F359").Value = 321196
$ws.Range("G359").Value = 3345
$ws.Range("F360").Value = 750433
$ws.Range("G360").Value = 5136
$ws.Range("F361").Value = 332579
$ws.Range("G361").Value = 2622

# Append new row 401 for 2021-04-10 (so 10. 04. 2021)
$ws.Range("A401").Value = 44295
$ws.Range("B401").Value = 370473
$ws.Range("C401").Value = 11753
$ws.Range("D401").Value = 1080
$ws.Range("E401").Value = 10487
$ws.Range("F401").Value = 944
$ws.Range("G401").Value = 2
